$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195, shifting existing rows 195:205 down to 196:206
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new weekly record
$ws.Range("A195").Value = 11
$ws.Range("B195").Value = "Vega Monumental Concepción"
$ws.Range("C195").Value = "Bíobío"
$ws.Range("D195").Value = 44826
$ws.Range("E195").Value = 8
$ws.Range("F195").Value = "Fruta"
$ws.Range("G195").Value = 100108
$ws.Range("H195").Value = "Tropicales y subtropicales"
$ws.Range("I195").Value = 100108005
$ws.Range("J195").Value = "Piña"
$ws.Range("K195").Value = "Sin especificar"
$ws.Range("L195").Value = "Tercera"
$ws.Range("M195").Value = 220
$ws.Range("N195").Value = 20000
$ws.Range("O195").Value = 21000
$ws.Range("P195").Value = 20545
$ws.Range("Q195").Value = "$/caja 16 unidades"
$ws.Range("R195").Value = "Ecuador"
$ws.Range("S195").Value = 1284
$ws.Range("T195").Value = 16
